$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing columns right
# (old A->B, B->C, C->D, D->E) and widths travel with their cells.
$ws.Columns("A:A").Insert()

# Updated Cypher query text for the StatQuery / query cells. Set these
# first (matches the shared-string ordering produced by the original
# edit), then fill in the new leading "TabName" / "CasesTab" column.
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$query = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$ws.Range("C2").Value = $statQuery
$ws.Range("B2").Value = $query

$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# New narrow first column, sized to fit "TabName"/"CasesTab".
$ws.Columns("A:A").ColumnWidth = 7.92

# Row 2 grew taller to fit the longer wrapped query text.
$ws.Rows("2:2").RowHeight = 174

# Selection moves off the data after the edits.
$ws.Range("B4").Select()

Write-Host "Edit applied"
